$wb = $excel.ActiveWorkbook

# --- Rename sheets (task order IDs refreshed) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16509961563752298"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1650996159753573"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1650996159753573"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16509961598096113"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16509961598736072"

# --- Sheet 1 (GNG) stimulus file updates ---
$ws1.Cells.Item(2, 2).Value = "go_stims-16509961563351939.csv"
$ws1.Cells.Item(3, 2).Value = "GNG_stims-16509961563592303.csv"
$ws1.Cells.Item(4, 2).Value = "go_stims-16509961563592303.csv"
$ws1.Cells.Item(5, 2).Value = "GNG_stims-16509961563752298.csv"

# --- Sheet 2 (NB) stimulus file updates ---
$ws2.Cells.Item(2, 2).Value = "ZB-match_5-16509961564951947.csv"
$ws2.Cells.Item(3, 2).Value = "ZB-match_2-16509961564552293.csv"
$ws2.Cells.Item(4, 2).Value = "OB-16509961572732677.csv"
$ws2.Cells.Item(5, 2).Value = "TB-16509961597376022.csv"
$ws2.Cells.Item(6, 2).Value = "ZB-match_1-16509961566492963.csv"
$ws2.Cells.Item(7, 2).Value = "OB-16509961576172962.csv"
$ws2.Cells.Item(8, 2).Value = "TB-1650996157737259.csv"
$ws2.Cells.Item(9, 2).Value = "OB-16509961573132637.csv"
$ws2.Cells.Item(10, 2).Value = "TB-16509961595696118.csv"

# --- Sheet 4 (TOL) stimulus file updates ---
$ws4.Cells.Item(2, 2).Value = "MM_stims-16509961597775793.csv"
$ws4.Cells.Item(3, 2).Value = "ZM_stims-1650996159753573.csv"
$ws4.Cells.Item(4, 2).Value = "MM_stims-1650996159793608.csv"
$ws4.Cells.Item(5, 2).Value = "ZM_stims-16509961597775793.csv"
$ws4.Cells.Item(6, 2).Value = "MM_stims-16509961598096113.csv"
$ws4.Cells.Item(7, 2).Value = "ZM_stims-1650996159793608.csv"

# --- Sheet 5 (vSAT) stimulus file updates ---
$ws5.Cells.Item(2, 2).Value = "vSAT_stims-1650996159857575.csv"
$ws5.Cells.Item(3, 2).Value = "SAT_stims-16509961598256018.csv"
$ws5.Cells.Item(4, 2).Value = "SAT_stims-16509961598096113.csv"
$ws5.Cells.Item(5, 2).Value = "vSAT_stims-16509961598415728.csv"
